$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 3rd bullet ("Create a EatGhost state.") - remove its text and its
# pre-existing "_GoBack" bookmark, leaving an empty list paragraph.
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$p5 = $d.Paragraphs.Item(5)
$p5Body = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$p5Body.Text = ""

# ------------------------------------------------------------------
# 2nd bullet - replace the sentence text, then attach a fresh
# "_GoBack" bookmark (collapsed, right after the new text).
# ------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Find.Execute( `
    "On next level, make sure that frightened mode is reset and that ghost return to scatter mode.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Make sure that ghosts are is switching between scatter and chase modes.", 2)

$p4 = $d.Paragraphs.Item(4)
$endPos = $p4.Range.End - 1

# Placing a zero-length bookmark exactly at the paragraph-end boundary
# is unreliable, so temporarily pad with a placeholder character,
# anchor the bookmark before it, then remove the placeholder again.
$d.Range($endPos, $endPos).InsertAfter("#")
$d.Bookmarks.Add("_GoBack", $d.Range($endPos, $endPos))
$d.Range($endPos, $endPos + 1).Text = ""

# ------------------------------------------------------------------
# 1st bullet - collapse the multi-run sentence (with proofErr marks)
# into a single new sentence.
# ------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Find.Execute( `
    "Make sure that going from frightened mode back to normal and back to frightened mode in the same level transitions without error.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Have Red ghost moving again.", 2)
